# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 260; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
